# Auto-generated edit script: refresh Tonberry Profits market data values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (ALC)
$ws.Range("H17").Value = 1447.2727
$ws.Range("J17").Value = 1447.2727
$ws.Range("L17").Value = 4341.8181
$ws.Range("N17").Value = -4677.8181

# Row 19 (ALC)
$ws.Range("H19").Value = 1357.5769
$ws.Range("J19").Value = 1484.7693
$ws.Range("L19").Value = 1484.7693
$ws.Range("N19").Value = -1834.7693

# Row 62 (ALC)
$ws.Range("H62").Value = 4665.6665
$ws.Range("J62").Value = 7198.2
$ws.Range("L62").Value = 7198.2
$ws.Range("N62").Value = -8446.200000000001

# Row 65 (ALC)
$ws.Range("H65").Value = 4665.6665
$ws.Range("J65").Value = 7198.2
$ws.Range("L65").Value = 35991
$ws.Range("N65").Value = -42231

# Row 106 (ALC)
$ws.Range("H106").Value = 1251.1818
$ws.Range("I106").Value = 1251.1818
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 1251.1818
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -620.1818000000001
$ws.Range("N106").ClearContents()

# Row 107 (ALC)
$ws.Range("H107").Value = 711.8570999999999
$ws.Range("I107").Value = 498.66666
$ws.Range("K107").Value = 498.66666
$ws.Range("M107").Value = 1421.33334

# Row 129 (ALC)
$ws.Range("H129").Value = 962.0909
$ws.Range("J129").Value = 929.5925999999999
$ws.Range("L129").Value = 2788.7778
$ws.Range("N129").Value = -12788.7778

# Row 137 (ALC)
$ws.Range("H137").Value = 1796.4193
$ws.Range("J137").Value = 1928.3462
$ws.Range("L137").Value = 5785.0386
$ws.Range("N137").Value = -10885.0386

# Row 138 (ALC)
$ws.Range("H138").Value = 1848.6111
$ws.Range("I138").Value = 1434.4546
$ws.Range("J138").Value = 2499.4285
$ws.Range("K138").Value = 4303.3638
$ws.Range("L138").Value = 7498.2855
$ws.Range("M138").Value = 836.6361999999999
$ws.Range("N138").Value = -17778.2855

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 3061.2615
$ws.Range("I32").Value = 2075.7817
$ws.Range("J32").Value = 8481.4
$ws.Range("K32").Value = 2075.7817
$ws.Range("L32").Value = 8481.4
$ws.Range("M32").Value = -1788.7817
$ws.Range("N32").Value = -9055.4

# Row 61 (ARM)
$ws.Range("H61").Value = 3239.2222
$ws.Range("I61").Value = 2230.182
$ws.Range("K61").Value = 2230.182
$ws.Range("M61").Value = -2018.182

# Row 74 (ARM)
$ws.Range("H74").Value = 1269.2174
$ws.Range("I74").Value = 565.2
$ws.Range("J74").Value = 2589.25
$ws.Range("K74").Value = 565.2
$ws.Range("L74").Value = 2589.25
$ws.Range("M74").Value = 308.8
$ws.Range("N74").Value = -4337.25

# Row 77 (ARM)
$ws.Range("H77").Value = 1269.2174
$ws.Range("I77").Value = 565.2
$ws.Range("J77").Value = 2589.25
$ws.Range("K77").Value = 2826
$ws.Range("L77").Value = 12946.25
$ws.Range("M77").Value = 1542
$ws.Range("N77").Value = -21682.25

# Row 132 (ARM)
$ws.Range("H132").Value = 1967.5366
$ws.Range("I132").Value = 1713.0834
$ws.Range("K132").Value = 5139.2502
$ws.Range("M132").Value = -2609.2502

# Row 136 (ARM)
$ws.Range("H136").Value = 3239.2222
$ws.Range("I136").Value = 2230.182
$ws.Range("K136").Value = 6690.545999999999
$ws.Range("M136").Value = -4140.545999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 122 (BSM)
$ws.Range("H122").Value = 68000
$ws.Range("J122").Value = 68000
$ws.Range("L122").Value = 68000
$ws.Range("N122").Value = -77800

# Row 134 (BSM)
$ws.Range("H134").Value = 5564.8887
$ws.Range("I134").Value = 5802.1665
$ws.Range("K134").Value = 17406.4995
$ws.Range("M134").Value = -14871.4995

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 1497
$ws.Range("I31").Value = 994.5714
$ws.Range("J31").Value = 1703.8823
$ws.Range("K31").Value = 994.5714
$ws.Range("L31").Value = 1703.8823
$ws.Range("M31").Value = -699.5714
$ws.Range("N31").Value = -2293.8823

# Row 34 (CRP)
$ws.Range("H34").Value = 1497
$ws.Range("I34").Value = 994.5714
$ws.Range("J34").Value = 1703.8823
$ws.Range("K34").Value = 994.5714
$ws.Range("L34").Value = 1703.8823
$ws.Range("M34").Value = -792.5714
$ws.Range("N34").Value = -2107.8823

# Row 94 (CRP)
$ws.Range("H94").Value = 1282.8
$ws.Range("I94").Value = 1200
$ws.Range("J94").Value = 1338
$ws.Range("K94").Value = 1200
$ws.Range("L94").Value = 1338
$ws.Range("M94").Value = -749
$ws.Range("N94").Value = -2240

# Row 132 (CRP)
$ws.Range("H132").Value = 3529.0833
$ws.Range("I132").Value = 2597.75
$ws.Range("K132").Value = 7793.25
$ws.Range("M132").Value = -5263.25

$ws = $wb.Worksheets.Item("CUL")
# Row 107 (CUL)
$ws.Range("H107").Value = 1650.6666
$ws.Range("J107").Value = 2341.7
$ws.Range("L107").Value = 7025.099999999999
$ws.Range("N107").Value = -10865.1

# Row 131 (CUL)
$ws.Range("H131").Value = 5162875
$ws.Range("J131").Value = 8758.088
$ws.Range("L131").Value = 26274.264
$ws.Range("N131").Value = -36354.264

$ws = $wb.Worksheets.Item("GSM")
# Row 23 (GSM)
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()

# Row 126 (GSM)
$ws.Range("H126").Value = 1618702.1
$ws.Range("I126").Value = 1918629.8
$ws.Range("K126").Value = 5755889.4
$ws.Range("M126").Value = -5753419.4

# Row 132 (GSM)
$ws.Range("H132").Value = 1925739.2
$ws.Range("I132").Value = 2566132.8
$ws.Range("J132").Value = 4559.2
$ws.Range("K132").Value = 7698398.399999999
$ws.Range("L132").Value = 13677.6
$ws.Range("M132").Value = -7695868.399999999
$ws.Range("N132").Value = -18737.6

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 4241.467
$ws.Range("I7").Value = 5001.3335
$ws.Range("J7").Value = 4051.5
$ws.Range("K7").Value = 5001.3335
$ws.Range("L7").Value = 4051.5
$ws.Range("M7").Value = -4889.3335
$ws.Range("N7").Value = -4275.5

# Row 40 (LTW)
$ws.Range("H40").Value = 8085.357
$ws.Range("I40").Value = 4114
$ws.Range("J40").Value = 12056.714
$ws.Range("K40").Value = 4114
$ws.Range("L40").Value = 12056.714
$ws.Range("M40").Value = -3978
$ws.Range("N40").Value = -12328.714

# Row 61 (LTW)
$ws.Range("H61").Value = 3498.8462
$ws.Range("I61").Value = 2950
$ws.Range("J61").Value = 3742.7778
$ws.Range("K61").Value = 2950
$ws.Range("L61").Value = 3742.7778
$ws.Range("M61").Value = -2748
$ws.Range("N61").Value = -4146.7778

# Row 113 (LTW)
$ws.Range("H113").Value = 3498.8462
$ws.Range("I113").Value = 2950
$ws.Range("J113").Value = 3742.7778
$ws.Range("L113").Value = 3742.7778
$ws.Range("M113").Value = -780
$ws.Range("N113").Value = -8082.7778

# Row 122 (LTW)
$ws.Range("H122").Value = 7191.3335
$ws.Range("I122").Value = 6419.3335
$ws.Range("J122").Value = 8349.333000000001
$ws.Range("K122").Value = 19258.0005
$ws.Range("L122").Value = 25047.999
$ws.Range("M122").Value = -16808.0005
$ws.Range("N122").Value = -29947.999

# Row 126 (LTW)
$ws.Range("H126").Value = 4241.467
$ws.Range("I126").Value = 5001.3335
$ws.Range("J126").Value = 4051.5
$ws.Range("K126").Value = 15004.0005
$ws.Range("L126").Value = 12154.5
$ws.Range("M126").Value = -12534.0005
$ws.Range("N126").Value = -17094.5

# Row 132 (LTW)
$ws.Range("H132").Value = 4201.125
$ws.Range("I132").Value = 1570.3334
$ws.Range("J132").Value = 4576.952
$ws.Range("K132").Value = 4711.0002
$ws.Range("L132").Value = 13730.856
$ws.Range("M132").Value = -2181.0002
$ws.Range("N132").Value = -18790.856

# Row 136 (LTW)
$ws.Range("H136").Value = 3152.2424
$ws.Range("I136").Value = 1890.1666
$ws.Range("K136").Value = 5670.4998
$ws.Range("M136").Value = -3120.4998

$ws = $wb.Worksheets.Item("WVR")
# Row 126 (WVR)
$ws.Range("H126").Value = 14224
$ws.Range("I126").Value = 18683.285
$ws.Range("K126").Value = 56049.855
$ws.Range("M126").Value = -53579.855

# Row 132 (WVR)
$ws.Range("H132").Value = 1413.6792
$ws.Range("I132").Value = 1053.3334
$ws.Range("K132").Value = 3160.0002
$ws.Range("M132").Value = -630.0001999999999

# Row 136 (WVR)
$ws.Range("H136").Value = 10686778
$ws.Range("I136").Value = 17923994
$ws.Range("J136").Value = 3270.9524
$ws.Range("K136").Value = 53771982
$ws.Range("L136").Value = 9812.8572
$ws.Range("M136").Value = -53769432
